# Update countries & provincias Spain
# Applies the changes described by the commit: refreshed COVID-19 numbers
# for several countries, a re-sorted pair of countries (Aruba/Bahamas and
# Timor Oriental/Santa Lucia swap rows because of the updated figures),
# and a refreshed "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (A1) -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 20:05"

# --- Country name swap: Aruba (row 135) <-> Bahamas (row 136) -------------
$ws.Range("A135").Value = "Bahamas"
$ws.Range("A136").Value = "Aruba"

# --- Country name swap: Timor Oriental (row 207) <-> Santa Lucia (row 208) -
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Timor Oriental"

# --- Updated statistics -----------------------------------------------------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7299820
$ws.Range("C4").Value = 12259
$ws.Range("D4").Value = 4536625
$ws.Range("E4").Value = 2553919
$ws.Range("G4").Value = 99
$ws.Range("H4").Value = 209276

# Row 5 - India
$ws.Range("B5").Value = 6050875
$ws.Range("C5").Value = 60294
$ws.Range("D5").Value = 4993558
$ws.Range("E5").Value = 962216
$ws.Range("G5").Value = 567
$ws.Range("H5").Value = 95101

# Row 14 - Francia
$ws.Range("B14").Value = 538569
$ws.Range("C14").Value = 11123
$ws.Range("E14").Value = 411951
$ws.Range("G14").Value = 27
$ws.Range("H14").Value = 31727

# Row 25 - Alemania
$ws.Range("B25").Value = 285821
$ws.Range("C25").Value = 796
$ws.Range("E25").Value = 26787

# Row 34 - Marruecos
$ws.Range("B34").Value = 117685
$ws.Range("C34").Value = 2444
$ws.Range("D34").Value = 95591
$ws.Range("E34").Value = 20025
$ws.Range("G34").Value = 28
$ws.Range("H34").Value = 2069

# Row 100 - Guinea
$ws.Range("B100").Value = 10548
$ws.Range("C100").Value = 36
$ws.Range("D100").Value = 9865
$ws.Range("E100").Value = 617
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 66

# Row 102 - Maldivas
$ws.Range("B102").Value = 10098
$ws.Range("C102").Value = 53
$ws.Range("D102").Value = 8847
$ws.Range("E102").Value = 1217

# Row 116 - Malaui
$ws.Range("B116").Value = 5768
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 4206
$ws.Range("E116").Value = 1383

# Row 135 - now Bahamas (previously Aruba's row)
$ws.Range("B135").Value = 3838
$ws.Range("C135").Value = 48
$ws.Range("D135").Value = 2005
$ws.Range("E135").Value = 1744
$ws.Range("H135").Value = 89

# Row 136 - now Aruba (previously Bahamas' row)
$ws.Range("B136").Value = 3832
$ws.Range("D136").Value = 2829
$ws.Range("E136").Value = 978
$ws.Range("H136").Value = 25

# Row 152 - Sierra Leona
$ws.Range("B152").Value = 2215
$ws.Range("C152").Value = 7
$ws.Range("D152").Value = 1681
$ws.Range("E152").Value = 462

Write-Output "Applied country/provincias updates"
